$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 update: convert phone number from text "+919876543210" to numeric 919876543210 ---
$ws.Range("C2").Value = 919876543210

# --- Row 3: new call record ---
$ws.Range("A3").Value = "call_20250703_102913_fa3bd916"
$ws.Range("B3").Value = "Vanshika panjwani"
$ws.Range("C3").Value = 917823844614
$ws.Range("D3").Value = "2025-07-03 04:59:13 IST"
$ws.Range("E3").Value = "0:01:34"
$ws.Range("F3").Value = "follow_up_needed"
$ws.Range("G3").Value = 23

# --- Row 4: new call record ---
$ws.Range("A4").Value = "call_20250703_103052_a13f7173"
$ws.Range("B4").Value = "Unknown"
$ws.Range("C4").Value = "Unknown"
$ws.Range("D4").Value = "2025-07-03 05:00:52 IST"
$ws.Range("E4").Value = "0:00:39"
$ws.Range("F4").Value = "call_incomplete"
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = "This call transcript records an unsuccessful attempt by an IVF clinic to connect with a potential patient.`n**1. Call Purpose and Context:**`n*   The call was initiated by Ritika from Aveya IVF – Rajouri Garden.`n*   The purpose was to follow up on a form recently submitted by the unknown recipient (or a family member) requesting `"clarity regarding fertility.`" This indicates a pre-existing inquiry or lead.`n**2. Patient's Main Concerns:**`n*   Based on the submitted form, the patient's primary concern was a general need for `"clarity on fertility`" issues. No specific details about the nature of this concern (e.g., difficulty conceiving, exploring options) were discussed.`n**3. Relevant Medical History Mentioned:**`n*   None. The call did not progress to a point where any medical history could be discussed.`n**4. Current Fertility Status:**`n*   Unknown. No information regarding the patient's current fertility status was exchanged.`n**5. Any Appointment Details Discussed:**`n*   None. Due to the call's repeated disconnections, no appointment details were discussed or scheduled.`n**6. Final Outcome of the Call:**`n*   The call was unsuccessful due to persistent technical issues leading to multiple disconnections.`n*   Ritika, the clinic representative, concluded the call, inviting the recipient to contact the clinic again at their convenience."
